# Apply edits to the "TestData" worksheet: update postcodes for existing
# AddCustomerTest rows, add a new AddCustomerTest row (Vasya Vasiliev),
# rename the customer values for OpenAccountTest, and add a new
# OpenAccountTest row (Ron Weasly / Pound).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")
$ws.Activate()

# Clear everything from row 5 down, since the OpenAccountTest block shifts
# down by one row and gains a new row; rebuild it cleanly to avoid leaving
# stale cell values behind.
$ws.Range("A5:E11").ClearContents()

# --- AddCustomerTest block ---
# Row 3 (Ivan Ivanov): postcode e3r4t5 -> dhfgh
$ws.Range("D3").Value = "dhfgh"

# Row 4 (Petr Petrov): postcode 2af4g5 -> fghfhf
$ws.Range("D4").Value = "fghfhf"

# New row 5: Vasya Vasiliev
$ws.Range("A5").Value = "Y"
$ws.Range("B5").Value = "Vasya"
$ws.Range("C5").Value = "Vasiliev"
$ws.Range("D5").Value = "fghfhf"
$ws.Range("E5").Value = "chrome"

# --- OpenAccountTest block (shifted down by one row, now starting at row 7) ---
$ws.Range("A7").Value = "OpenAccountTest"

$ws.Range("A8").Value = "Runmode"
$ws.Range("B8").Value = "customer"
$ws.Range("C8").Value = "currency"
$ws.Range("D8").Value = "browser"

$ws.Range("A9").Value = "Y"
$ws.Range("B9").Value = "Hermoine Granger"
$ws.Range("C9").Value = "Dollar"
$ws.Range("D9").Value = "chrome"

$ws.Range("A10").Value = "Y"
$ws.Range("B10").Value = "Harry Potter"
$ws.Range("C10").Value = "Rupee"
$ws.Range("D10").Value = "firefox"

# New row 11: Ron Weasly
$ws.Range("A11").Value = "Y"
$ws.Range("B11").Value = "Ron Weasly"
$ws.Range("C11").Value = "Pound"
$ws.Range("D11").Value = "firefox"

# Update the selected cell to match the new state.
$ws.Range("D6").Select()
